$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 36

$data = @(
    @("Edery", "Aaron", "Oui", 4),
    @("Bitton", "Dan", "Oui", 2),
    @("Bitton", "Dan", "Oui", 2),
    @("Bitton", "Dan", "Oui", 2),
    @("Bitton", "Dan", "Oui", 2),
    @("Bitton", "Dan", "Oui", 2),
    @("Bitton", "Dan", "Oui", 2),
    @("Edery", "Dan", "Oui", 4),
    @("yaire", "coco", "Oui", 3),
    @("malik", "jojo", "Oui", 1)
)

$row = 8
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $row++
}
